# matrix.xlsx edit: refresh the sample matrix values, drop the stray
# number-format style that was on A2, and move the selection cursor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 no longer needs the "#,##0" number-format style (cellXfs index 1) -
# put it back on the default/Normal style so that xf record becomes unused.
$ws.Cells.Item(2, 1).Style = "Normal"

# Refresh the small 3x2 block of sample numbers.
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = 10
$ws.Range("A3").Value = 20
$ws.Range("B3").Value = 20
$ws.Range("C3").Value = 10

# Move the active selection to C8 (was E2).
$ws.Range("C8").Select() | Out-Null
